$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "bateria"
$ws.Range("B11").Value = 500
$ws.Range("A12").Value = "bomba de agua"
$ws.Range("B13").Value = 400

$ws.Range("A9").Select()
